$d = $word.ActiveDocument

$replacements = @(
    @("91×19=1729", "50×63=3150"),
    @("96×73=7008", "87×53=4611"),
    @("38×28=1064", "70×43=3010"),
    @("24×13=312",  "91×79=7189"),
    @("72×47=3384", "78×25=1950"),
    @("70×45=3150", "80×11=880"),
    @("34×80=2720", "16×45=720"),
    @("87×84=7308", "98×81=7938"),
    @("96×69=6624", "43×53=2279"),
    @("20×58=1160", "22×19=418"),
    @("76×84=6384", "24×20=480"),
    @("45×63=2835", "47×53=2491"),
    @("96×58=5568", "85×86=7310"),
    @("64×85=5440", "25×29=725"),
    @("34×63=2142", "69×95=6555"),
    @("57×54=3078", "59×33=1947"),
    @("84×22=1848", "93×33=3069"),
    @("37×22=814",  "87×94=8178"),
    @("85×93=7905", "35×98=3430"),
    @("80×69=5520", "49×27=1323"),
    @("91×72=6552", "70×99=6930"),
    @("72×50=3600", "52×28=1456"),
    @("45×28=1260", "67×36=2412"),
    @("34×73=2482", "81×33=2673"),
    @("86×73=6278", "58×33=1914")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
